$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2025-04-14 10:00:00"
$ws.Range("B2").Value = "admin"
$ws.Range("C2").Value = "inventario"
$ws.Range("D2").Value = "inserción"
$ws.Range("E2").Value = "Agregó un nuevo ítem"
$ws.Range("F2").Value = "192.168.1.1"

$ws.Range("A3").Value = "2025-04-14 11:00:00"
$ws.Range("B3").Value = "user1"
$ws.Range("C3").Value = "logística"
$ws.Range("D3").Value = "modificación"
$ws.Range("E3").Value = "Actualizó estado de entrega"
$ws.Range("F3").Value = "192.168.1.2"
